$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Mapping of (row, col) -> new text. Row numbers are the Word table row
# indices (1-based) that actually contain data (1, 5, 9, 13, 17); the
# rows in between are empty spacer rows and are left untouched.

$updates = @{
    # Data row 1 (source row index 1)
    "1,1" = "85÷5=17, 0"
    "1,2" = "24÷3=8, 0"
    "1,3" = "38÷2=19, 0"
    "1,4" = "61÷3=20, 1"
    "1,5" = "94÷3=31, 1"

    # Data row 2 (source row index 5)
    "5,1" = "12÷6=2, 0"
    "5,2" = "98÷7=14, 0"
    "5,3" = "47÷3=15, 2"
    "5,4" = "32÷7=4, 4"
    "5,5" = "81÷4=20, 1"

    # Data row 3 (source row index 9)
    "9,1" = "36÷4=9, 0"
    "9,2" = "42÷8=5, 2"
    "9,3" = "83÷3=27, 2"
    "9,4" = "66÷7=9, 3"
    "9,5" = "10÷4=2, 2"

    # Data row 4 (source row index 13)
    "13,1" = "50÷7=7, 1"
    "13,2" = "90÷6=15, 0"
    "13,3" = "67÷7=9, 4"
    "13,4" = "26÷9=2, 8"
    "13,5" = "24÷8=3, 0"

    # Data row 5 (source row index 17)
    "17,1" = "68÷8=8, 4"
    "17,2" = "42÷7=6, 0"
    "17,3" = "19÷3=6, 1"
    "17,4" = "92÷7=13, 1"
    "17,5" = "91÷4=22, 3"
}

foreach ($key in $updates.Keys) {
    $parts = $key.Split(",")
    $row = [int]$parts[0]
    $col = [int]$parts[1]
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $updates[$key]
}
